$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 3294.5
$ws.Cells.Item(100, 9).Value = 1590
$ws.Cells.Item(100, 11).Value = 1590
$ws.Cells.Item(100, 13).Value = -1049
$ws.Cells.Item(132, 8).Value = 3276.4375
$ws.Cells.Item(132, 9).Value = 3304.8667
$ws.Cells.Item(132, 11).Value = 9914.6001
$ws.Cells.Item(132, 13).Value = -7384.6001
$ws.Cells.Item(137, 8).Value = 3552.1628
$ws.Cells.Item(137, 9).Value = 3143.7334
$ws.Cells.Item(137, 10).Value = 4494.6924
$ws.Cells.Item(137, 11).Value = 9431.200199999999
$ws.Cells.Item(137, 12).Value = 13484.0772
$ws.Cells.Item(137, 13).Value = -6881.200199999999
$ws.Cells.Item(137, 14).Value = -18584.0772
$ws.Cells.Item(138, 8).Value = 4631.151
$ws.Cells.Item(138, 10).Value = 4500.9424
$ws.Cells.Item(138, 12).Value = 13502.8272
$ws.Cells.Item(138, 14).Value = -23782.8272

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 5665.5
$ws.Cells.Item(74, 9).Value = 3062.5217
$ws.Cells.Item(74, 10).Value = 14218.143
$ws.Cells.Item(74, 11).Value = 3062.5217
$ws.Cells.Item(74, 12).Value = 14218.143
$ws.Cells.Item(74, 13).Value = -2188.5217
$ws.Cells.Item(74, 14).Value = -15966.143
$ws.Cells.Item(77, 8).Value = 5665.5
$ws.Cells.Item(77, 9).Value = 3062.5217
$ws.Cells.Item(77, 10).Value = 14218.143
$ws.Cells.Item(77, 11).Value = 15312.6085
$ws.Cells.Item(77, 12).Value = 71090.715
$ws.Cells.Item(77, 13).Value = -10944.6085
$ws.Cells.Item(77, 14).Value = -79826.715
$ws.Cells.Item(120, 8).Value = 65777
$ws.Cells.Item(120, 10).Value = 65777
$ws.Cells.Item(120, 12).Value = 65777
$ws.Cells.Item(120, 14).Value = -75453
$ws.Cells.Item(125, 8).Value = 96000
$ws.Cells.Item(125, 10).Value = 96000
$ws.Cells.Item(125, 12).Value = 96000
$ws.Cells.Item(125, 14).Value = -105840
$ws.Cells.Item(132, 8).Value = 7733.6665
$ws.Cells.Item(132, 9).Value = 2810.1667
$ws.Cells.Item(132, 10).Value = 17580.666
$ws.Cells.Item(132, 11).Value = 8430.500100000001
$ws.Cells.Item(132, 12).Value = 52741.99800000001
$ws.Cells.Item(132, 13).Value = -5900.500100000001
$ws.Cells.Item(132, 14).Value = -57801.99800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1959.1666
$ws.Cells.Item(16, 9).Value = 1031.8572
$ws.Cells.Item(16, 10).Value = 3257.4
$ws.Cells.Item(16, 11).Value = 1031.8572
$ws.Cells.Item(16, 12).Value = 3257.4
$ws.Cells.Item(16, 13).Value = -744.8571999999999
$ws.Cells.Item(16, 14).Value = -3831.4
$ws.Cells.Item(31, 8).Value = 8505.023999999999
$ws.Cells.Item(31, 9).Value = 4439.7085
$ws.Cells.Item(31, 11).Value = 4439.7085
$ws.Cells.Item(31, 13).Value = -4144.7085
$ws.Cells.Item(34, 8).Value = 8505.023999999999
$ws.Cells.Item(34, 9).Value = 4439.7085
$ws.Cells.Item(34, 11).Value = 4439.7085
$ws.Cells.Item(34, 13).Value = -4237.7085
$ws.Cells.Item(58, 8).Value = 17999.5
$ws.Cells.Item(58, 10).Value = 17999.5
$ws.Cells.Item(58, 12).Value = 17999.5
$ws.Cells.Item(58, 14).Value = -18405.5
$ws.Cells.Item(99, 8).Value = 30099.75
$ws.Cells.Item(99, 9).Value = 37233.332
$ws.Cells.Item(99, 10).Value = 8699
$ws.Cells.Item(99, 11).Value = 37233.332
$ws.Cells.Item(99, 12).Value = 8699
$ws.Cells.Item(99, 13).Value = -35735.332
$ws.Cells.Item(99, 14).Value = -11695
$ws.Cells.Item(113, 8).Value = 1959.1666
$ws.Cells.Item(113, 9).Value = 1031.8572
$ws.Cells.Item(113, 10).Value = 3257.4
$ws.Cells.Item(113, 11).Value = 1031.8572
$ws.Cells.Item(113, 12).Value = 3257.4
$ws.Cells.Item(113, 13).Value = 1138.1428
$ws.Cells.Item(113, 14).Value = -7597.4
$ws.Cells.Item(122, 8).Value = 3998.805
$ws.Cells.Item(122, 9).Value = 3730.2593
$ws.Cells.Item(122, 10).Value = 4516.7144
$ws.Cells.Item(122, 11).Value = 11190.7779
$ws.Cells.Item(122, 12).Value = 13550.1432
$ws.Cells.Item(122, 13).Value = -8740.777900000001
$ws.Cells.Item(122, 14).Value = -18450.1432
$ws.Cells.Item(126, 8).Value = 30099.75
$ws.Cells.Item(126, 9).Value = 37233.332
$ws.Cells.Item(126, 10).Value = 8699
$ws.Cells.Item(126, 11).Value = 111699.996
$ws.Cells.Item(126, 12).Value = 26097
$ws.Cells.Item(126, 13).Value = -109229.996
$ws.Cells.Item(126, 14).Value = -31037
$ws.Cells.Item(132, 8).Value = 8895.647000000001
$ws.Cells.Item(132, 9).Value = 7373.357
$ws.Cells.Item(132, 10).Value = 15999.667
$ws.Cells.Item(132, 11).Value = 22120.071
$ws.Cells.Item(132, 12).Value = 47999.001
$ws.Cells.Item(132, 13).Value = -19590.071
$ws.Cells.Item(132, 14).Value = -53059.001
$ws.Cells.Item(134, 8).Value = 3747.4856
$ws.Cells.Item(134, 9).Value = 3276.2964
$ws.Cells.Item(134, 10).Value = 5337.75
$ws.Cells.Item(134, 11).Value = 9828.889200000001
$ws.Cells.Item(134, 12).Value = 16013.25
$ws.Cells.Item(134, 13).Value = -7293.889200000001
$ws.Cells.Item(134, 14).Value = -21083.25
$ws.Cells.Item(136, 8).Value = 17999.5
$ws.Cells.Item(136, 10).Value = 17999.5
$ws.Cells.Item(136, 12).Value = 53998.5
$ws.Cells.Item(136, 14).Value = -59098.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1680
$ws.Cells.Item(80, 9).Value = 1620
$ws.Cells.Item(80, 10).Value = 1800
$ws.Cells.Item(80, 11).Value = 1620
$ws.Cells.Item(80, 12).Value = 1800
$ws.Cells.Item(80, 13).Value = -622
$ws.Cells.Item(80, 14).Value = -3796
$ws.Cells.Item(83, 8).Value = 1680
$ws.Cells.Item(83, 9).Value = 1620
$ws.Cells.Item(83, 10).Value = 1800
$ws.Cells.Item(83, 11).Value = 8100
$ws.Cells.Item(83, 12).Value = 9000
$ws.Cells.Item(83, 13).Value = -3108
$ws.Cells.Item(83, 14).Value = -18984
$ws.Cells.Item(102, 8).Value = 4247.1763
$ws.Cells.Item(102, 9).Value = 3243.7856
$ws.Cells.Item(102, 11).Value = 3243.7856
$ws.Cells.Item(102, 13).Value = -1621.7856
$ws.Cells.Item(122, 8).Value = 7023.36
$ws.Cells.Item(122, 9).Value = 4333.6875
$ws.Cells.Item(122, 11).Value = 13001.0625
$ws.Cells.Item(122, 13).Value = -10551.0625
$ws.Cells.Item(132, 8).Value = 5637.3555
$ws.Cells.Item(132, 9).Value = 5704.8647
$ws.Cells.Item(132, 10).Value = 5325.125
$ws.Cells.Item(132, 11).Value = 17114.5941
$ws.Cells.Item(132, 12).Value = 15975.375
$ws.Cells.Item(132, 13).Value = -14584.5941
$ws.Cells.Item(132, 14).Value = -21035.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 170166.33
$ws.Cells.Item(7, 9).Value = 253000
$ws.Cells.Item(7, 11).Value = 253000
$ws.Cells.Item(7, 13).Value = -252888
$ws.Cells.Item(22, 8).Value = 2214.1538
$ws.Cells.Item(22, 10).Value = 2833.3333
$ws.Cells.Item(22, 12).Value = 2833.3333
$ws.Cells.Item(22, 14).Value = -3423.3333
$ws.Cells.Item(27, 8).Value = 2214.1538
$ws.Cells.Item(27, 10).Value = 2833.3333
$ws.Cells.Item(27, 12).Value = 2833.3333
$ws.Cells.Item(27, 14).Value = -3047.3333
$ws.Cells.Item(40, 8).Value = 4955.6523
$ws.Cells.Item(40, 9).Value = 3735.7896
$ws.Cells.Item(40, 11).Value = 3735.7896
$ws.Cells.Item(40, 13).Value = -3599.7896
$ws.Cells.Item(82, 8).Value = 4604.636
$ws.Cells.Item(82, 9).Value = 4507.2856
$ws.Cells.Item(82, 11).Value = 4507.2856
$ws.Cells.Item(82, 13).Value = -4146.2856
$ws.Cells.Item(85, 8).Value = 4604.636
$ws.Cells.Item(85, 9).Value = 4507.2856
$ws.Cells.Item(85, 11).Value = 4507.2856
$ws.Cells.Item(85, 13).Value = -3259.2856
$ws.Cells.Item(122, 8).Value = 6129.2354
$ws.Cells.Item(122, 9).Value = 5190.636
$ws.Cells.Item(122, 11).Value = 15571.908
$ws.Cells.Item(122, 13).Value = -13121.908
$ws.Cells.Item(126, 8).Value = 170166.33
$ws.Cells.Item(126, 9).Value = 253000
$ws.Cells.Item(126, 11).Value = 759000
$ws.Cells.Item(126, 13).Value = -756530

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1300
$ws.Cells.Item(81, 9).Value = 1300
$ws.Cells.Item(81, 11).Value = 2600
$ws.Cells.Item(81, 13).Value = -1539
$ws.Cells.Item(84, 8).Value = 1300
$ws.Cells.Item(84, 9).Value = 1300
$ws.Cells.Item(84, 11).Value = 13000
$ws.Cells.Item(84, 13).Value = -7696
$ws.Cells.Item(100, 8).Value = 1296.5
$ws.Cells.Item(100, 9).Value = 809.3333
$ws.Cells.Item(100, 10).Value = 1783.6666
$ws.Cells.Item(100, 11).Value = 1618.6666
$ws.Cells.Item(100, 12).Value = 3567.3332
$ws.Cells.Item(100, 13).Value = -1077.6666
$ws.Cells.Item(100, 14).Value = -4649.3332
